# Updates cryptos list price/volume figures (and two coin-row swaps)
# to match the refreshed scrape, per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value looks like a plain decimal number
# (e.g. "680.65") must be entered with a leading apostrophe via .Formula
# so Excel keeps them as Text (matching the original inline-string cells)
# instead of silently re-typing them as Number (which would also drop
# trailing zeros, e.g. turning "1.00" into 1).

$ws.Range("D2").Value = '69.214.56'
$ws.Range("E2").Value = '  -2.11%  '
$ws.Range("D3").Value = '3.691.48'
$ws.Range("E3").Value = '  -2.85%  '
$ws.Range("D5").Formula = '''680.65'
$ws.Range("E5").Value = '  -3.66%  '
$ws.Range("D6").Formula = '''162.45'
$ws.Range("E6").Value = '  -4.43%  '
$ws.Range("D7").Value = '3.687.83'
$ws.Range("E7").Value = '  -2.92%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  -4.16%  '
$ws.Range("E10").Value = '  -7.20%  '
$ws.Range("D11").Formula = '''7.30'
$ws.Range("E11").Value = '  -1.06%  '
$ws.Range("E12").Value = '  -1.22%  '
$ws.Range("D13").Formula = '''0.0000237'
$ws.Range("E13").Value = '  -6.32%  '
$ws.Range("D14").Formula = '''33.59'
$ws.Range("E14").Value = '  -6.76%  '
$ws.Range("D15").Value = '4.315.49'
$ws.Range("E15").Value = '  -2.84%  '
$ws.Range("D16").Value = '3.691.39'
$ws.Range("E16").Value = '  -3.35%  '
$ws.Range("D17").Value = '69.306.58'
$ws.Range("E17").Value = '  -2.01%  '
$ws.Range("E18").Value = '  -1.76%  '
$ws.Range("D19").Formula = '''16.35'
$ws.Range("E19").Value = '  -5.95%  '
$ws.Range("D20").Formula = '''6.62'
$ws.Range("E20").Value = '  -6.96%  '
$ws.Range("D21").Formula = '''481.70'
$ws.Range("E21").Value = '  -3.09%  '
$ws.Range("D22").Formula = '''9.81'
$ws.Range("E22").Value = '  -7.39%  '
$ws.Range("D23").Formula = '''0.667'
$ws.Range("E23").Value = '  -8.32%  '
$ws.Range("D24").Formula = '''79.96'
$ws.Range("E24").Value = '  -5.42%  '
$ws.Range("D25").Value = '3.836.62'
$ws.Range("D26").Formula = '''0.0000129'
$ws.Range("E26").Value = '  -10.46%  '
$ws.Range("D27").Formula = '''11.52'
$ws.Range("E27").Value = '  -4.73%  '
$ws.Range("D28").Formula = '''1.00'
$ws.Range("E28").Value = '  -0.01%  '
$ws.Range("D29").Formula = '''9.59'
$ws.Range("E29").Value = '  -7.89%  '
$ws.Range("D30").Formula = '''1.84'
$ws.Range("E30").Value = '  -10.15%  '
$ws.Range("E31").Value = '  -10.64%  '
$ws.Range("D32").Formula = '''2.12'
$ws.Range("E32").Value = '  -4.47%  '
$ws.Range("D33").Formula = '''6.85'
$ws.Range("E33").Value = '  -6.16%  '
$ws.Range("D34").Formula = '''27.11'
$ws.Range("E34").Value = '  -6.48%  '
$ws.Range("B35").Value = 'Binance-PegBSC-USD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D35").Formula = '''1.00'
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").Formula = '''0.166'
$ws.Range("E36").Value = '  -3.86%  '
$ws.Range("D37").Value = '3.654.50'
$ws.Range("E37").Value = '  -3.09%  '
$ws.Range("D38").Formula = '''8.56'
$ws.Range("E38").Value = '  -5.85%  '
$ws.Range("D39").Formula = '''6.04'
$ws.Range("E39").Value = '  +1.96%  '
$ws.Range("D40").Formula = '''0.0943'
$ws.Range("E40").Value = '  -6.99%  '
$ws.Range("E41").Value = '  +0.02%  '
$ws.Range("D42").Formula = '''2.19'
$ws.Range("E42").Value = '  -5.79%  '
$ws.Range("E43").Value = '  -0.06%  '
$ws.Range("D44").Formula = '''0.961'
$ws.Range("E44").Value = '  -7.26%  '
$ws.Range("D45").Formula = '''158.86'
$ws.Range("E45").Value = '  -3.65%  '
$ws.Range("D46").Formula = '''48.15'
$ws.Range("E46").Value = '  -1.17%  '
$ws.Range("D47").Formula = '''2.85'
$ws.Range("E47").Value = '  -12.51%  '
$ws.Range("D48").Formula = '''0.000281'
$ws.Range("E48").Value = '  -12.81%  '
$ws.Range("B49").Value = 'ONDO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D49").Formula = '''1.30'
$ws.Range("E49").Value = '  -3.78%  '
$ws.Range("B50").Value = 'Bittensor'
$ws.Range("C50").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D50").Formula = '''388.16'
$ws.Range("E50").Value = '  -8.79%  '
$ws.Range("D51").Formula = '''8.10'
$ws.Range("E51").Value = '  -5.79%  '
